# Update "想去人数" (interest counts) for a few events in the
# "展览" and "全部类型" worksheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7257
$ws1.Range("F4").Value = 125
$ws1.Range("F7").Value = 96
$ws1.Range("F8").Value = 613

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7257
$ws4.Range("F5").Value = 125
$ws4.Range("F9").Value = 96
$ws4.Range("F10").Value = 613
